# 11.5.1.xlsx — add the "2020" column (Q) to the SDG 11.5.1 indicator table.
#
# The source table already holds one column per year (D:P = 2007..2019).
# This change appends a new year column, Q, for 2020, copying the number
# format/border/alignment from the corresponding 2019 cell (column P) in
# each row and then filling in the 2020 figures (or "-" where no data is
# available, reusing the shared "-" string already used elsewhere in the
# sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> 2020 value. $null marks the blank header-separator row (row 3,
# format only, no value) and "-" marks "no data" cells (shared string).
$values = [ordered]@{
    3  = $null
    4  = 2020
    5  = 51
    6  = 29
    7  = 22
    8  = 5
    9  = 3
    10 = 2
    11 = 15
    12 = 9
    13 = 5
    14 = "-"
    15 = "-"
    16 = "-"
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 7
    21 = 7
    22 = "-"
    23 = "-"
    24 = "-"
    25 = "-"
    26 = 24
    27 = 10
    28 = 14
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

foreach ($row in $values.Keys) {
    $src = $ws.Range("P$row")
    $dst = $ws.Range("Q$row")

    # Bring over the same number format / font / borders / alignment used
    # by the 2019 column for this row.
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats

    $v = $values[$row]
    if ($null -ne $v) {
        $dst.Value = $v
    }
}

$excel.CutCopyMode = $false

# Match the author's final selection.
[void]$ws.Range("H26").Select()
